# Update the "dSF" (column F) values on Sheet1 as per repull/recalculation of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = -7
$ws.Range("F3").Value = -5
$ws.Range("F5").Value = 5
$ws.Range("F6").Value = -4
$ws.Range("F7").Value = 3
